$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Andamento_2")
$ws.Range("C56").Value = 43921
$ws.Range("H56").Value = 404
$ws.Range("M56").Value = "S"
$ws.Range("B56:N56").Font.Color = 5287936

$ws.Range("C79").Value = 43922
$ws.Range("H79").Value = 30
$ws.Range("M79").Value = "S"
$ws.Range("B79:N79").Font.Color = 5287936

$ws.Range("C80").Value = 43922
$ws.Range("H80").Value = 56
$ws.Range("M80").Value = "S"
$ws.Range("B80:N80").Font.Color = 5287936
